# Update the dSF (column F) values on the active worksheet to reflect the
# re-pulled data / mean calculation described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = -1
    4  = 1
    7  = 0
    8  = 5
    9  = -1
    10 = -1
    11 = 6
    12 = -7
    13 = 1
    14 = 4
    15 = -3
    17 = 2
    18 = -4
    21 = 2
    22 = 2
    23 = -3
    24 = -2
    25 = 5
    26 = -5
    27 = -1
    28 = 1
    29 = 4
    30 = -1
    31 = 1
    32 = -1
    33 = -1
    35 = 0
    36 = -2
    39 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
